$d = $word.ActiveDocument

# Locate the paragraph that contains the certification sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*This is to certify*") {
        $target = $p
        break
    }
}

# Within that paragraph, find the "Let's Meet" project-name text so we can
# carve it out into its own (bold) run.
$nameRange = $target.Range.Duplicate
$nameRange.Find.ClearFormatting()
$nameRange.Find.Text = "Let's Meet"
$nameRange.Find.Execute() | Out-Null

# Wrap the project name in the (re-homed) "_GoBack" bookmark first - adding a
# bookmark with a name that already exists elsewhere in the document moves it
# here (deleting the old bookmarkStart/bookmarkEnd pair near "Date").
$d.Bookmarks.Add("_GoBack", $nameRange) | Out-Null

# Now make the project name bold; this splits the run the project name lives
# in away from the surrounding plain-formatted text on either side of it.
$nameRange.Font.Bold = 1
